# Applies the table-style change captured in the commit: the three tables
# that still carried the deck's default table style ("Table_0",
# {B35A0A91-E87D-4422-AB54-256141216A9B}) are switched to the built-in
# "No Style, No Grid" table style ({5CCA61C1-345E-4EFB-80F9-48927991862F}).

$p = $ppt.ActivePresentation
$noStyleNoGrid = "{5CCA61C1-345E-4EFB-80F9-48927991862F}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($noStyleNoGrid)
        }
    }
}
